# Auto-generated Excel COM-interop script to apply the Kujata_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for
# specific leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 24000
$ws.Range("J3").Value = 24000
$ws.Range("L3").Value = 24000
$ws.Range("N3").Value = -24228

$ws.Range("H33").Value = 503.15384
$ws.Range("I33").Value = 548.1
$ws.Range("K33").Value = 548.1
$ws.Range("M33").Value = -319.1

$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376

$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880

$ws.Range("H76").Value = 5736.2085
$ws.Range("I76").Value = 4955.143
$ws.Range("J76").Value = 6057.8237
$ws.Range("K76").Value = 4955.143
$ws.Range("L76").Value = 6057.8237
$ws.Range("M76").Value = -4640.143
$ws.Range("N76").Value = -6687.8237

$ws.Range("H79").Value = 5736.2085
$ws.Range("I79").Value = 4955.143
$ws.Range("J79").Value = 6057.8237
$ws.Range("K79").Value = 4955.143
$ws.Range("L79").Value = 6057.8237
$ws.Range("M79").Value = -3863.143
$ws.Range("N79").Value = -8241.823700000001

$ws.Range("H100").Value = 1943.8889
$ws.Range("I100").Value = 1811.875
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1811.875
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1270.875
$ws.Range("N100").Value = -4082

$ws.Range("H102").Value = 24000
$ws.Range("J102").Value = 24000
$ws.Range("L102").Value = 24000
$ws.Range("N102").Value = -30490

$ws.Range("H111").Value = 1016.3333
$ws.Range("I111").Value = 1026.2858
$ws.Range("K111").Value = 3078.8574
$ws.Range("M111").Value = -11.85740000000033

$ws.Range("H112").Value = 2260.889
$ws.Range("J112").Value = 2613.5173
$ws.Range("L112").Value = 7840.5519
$ws.Range("N112").Value = -10056.5519

$ws.Range("H125").Value = 1741.3334
$ws.Range("I125").Value = 1670.6666
$ws.Range("K125").Value = 15035.9994
$ws.Range("M125").Value = -12575.9994

$ws.Range("H132").Value = 11911093
$ws.Range("I132").Value = 13895662
$ws.Range("J132").Value = 3676.5
$ws.Range("K132").Value = 41686986
$ws.Range("L132").Value = 11029.5
$ws.Range("M132").Value = -41684456
$ws.Range("N132").Value = -16089.5

$ws.Range("H138").Value = 1480.0408
$ws.Range("J138").Value = 1893.2985
$ws.Range("L138").Value = 5679.895500000001
$ws.Range("N138").Value = -15959.8955

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1397.6666
$ws.Range("I61").Value = 885.2222
$ws.Range("K61").Value = 885.2222
$ws.Range("M61").Value = -673.2222

$ws.Range("H70").Value = 40000
$ws.Range("J70").Value = 40000
$ws.Range("L70").Value = 40000
$ws.Range("N70").Value = -40540

$ws.Range("H73").Value = 40000
$ws.Range("J73").Value = 40000
$ws.Range("L73").Value = 40000
$ws.Range("N73").Value = -41872

$ws.Range("H102").Value = 10419682
$ws.Range("I102").Value = 12823538
$ws.Range("J102").Value = 2970.3333
$ws.Range("K102").Value = 12823538
$ws.Range("L102").Value = 2970.3333
$ws.Range("M102").Value = -12821916
$ws.Range("N102").Value = -6214.3333

$ws.Range("H109").Value = 40500
$ws.Range("J109").Value = 40500
$ws.Range("L109").Value = 40500
$ws.Range("N109").Value = -43274

$ws.Range("H132").Value = 2244.0312
$ws.Range("I132").Value = 1912.68
$ws.Range("K132").Value = 5738.04
$ws.Range("M132").Value = -3208.04

$ws.Range("H136").Value = 1397.6666
$ws.Range("I136").Value = 885.2222
$ws.Range("K136").Value = 2655.6666
$ws.Range("M136").Value = -105.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 41668100
$ws.Range("I99").Value = 62501316
$ws.Range("K99").Value = 62501316
$ws.Range("M99").Value = -62499818

$ws.Range("H107").Value = 1823.2
$ws.Range("I107").Value = 1770.2222
$ws.Range("J107").Value = 2300
$ws.Range("K107").Value = 1770.2222
$ws.Range("L107").Value = 2300
$ws.Range("M107").Value = 149.7778000000001
$ws.Range("N107").Value = -6140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1196.2603
$ws.Range("I31").Value = 1069.8939
$ws.Range("J31").Value = 2387.7144
$ws.Range("K31").Value = 1069.8939
$ws.Range("L31").Value = 2387.7144
$ws.Range("M31").Value = -774.8939
$ws.Range("N31").Value = -2977.7144

$ws.Range("H34").Value = 1196.2603
$ws.Range("I34").Value = 1069.8939
$ws.Range("J34").Value = 2387.7144
$ws.Range("K34").Value = 1069.8939
$ws.Range("L34").Value = 2387.7144
$ws.Range("M34").Value = -867.8939
$ws.Range("N34").Value = -2791.7144

$ws.Range("H107").Value = 660.2222
$ws.Range("I107").Value = 450.5
$ws.Range("J107").Value = 828
$ws.Range("K107").Value = 450.5
$ws.Range("L107").Value = 828
$ws.Range("M107").Value = 1469.5
$ws.Range("N107").Value = -4668

$ws.Range("H132").Value = 2000.3103
$ws.Range("I132").Value = 1465.85
$ws.Range("J132").Value = 3188
$ws.Range("K132").Value = 4397.549999999999
$ws.Range("L132").Value = 9564
$ws.Range("M132").Value = -1867.549999999999
$ws.Range("N132").Value = -14624

$ws.Range("H134").Value = 1105.3667
$ws.Range("I134").Value = 933.7406999999999
$ws.Range("K134").Value = 2801.2221
$ws.Range("M134").Value = -266.2221

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 457.2143
$ws.Range("I7").Value = 479.9
$ws.Range("J7").Value = 400.5
$ws.Range("K7").Value = 1439.7
$ws.Range("L7").Value = 1201.5
$ws.Range("M7").Value = -1327.7
$ws.Range("N7").Value = -1425.5

$ws.Range("H34").Value = 1961.3462
$ws.Range("J34").Value = 2449.389
$ws.Range("L34").Value = 7348.167
$ws.Range("N34").Value = -7516.167

$ws.Range("H104").Value = 4190.2
$ws.Range("J104").Value = 5089.9
$ws.Range("L104").Value = 15269.7
$ws.Range("N104").Value = -20511.7

$ws.Range("H122").Value = 769.0769
$ws.Range("I122").Value = 446.6
$ws.Range("J122").Value = 970.625
$ws.Range("K122").Value = 4019.4
$ws.Range("L122").Value = 8735.625
$ws.Range("M122").Value = -1569.4

$ws.Range("H137").Value = 2105
$ws.Range("I137").Value = 795.38464
$ws.Range("K137").Value = 2386.15392
$ws.Range("M137").Value = 2713.84608

$ws.Range("H139").Value = 1635.075
$ws.Range("I139").Value = 1762.6818
$ws.Range("K139").Value = 5288.0454
$ws.Range("M139").Value = -148.0454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 43000
$ws.Range("J104").Value = 43000
$ws.Range("L104").Value = 43000
$ws.Range("N104").Value = -49988

$ws.Range("H107").Value = 783.9048
$ws.Range("I107").Value = 797.4375
$ws.Range("K107").Value = 797.4375
$ws.Range("M107").Value = 1122.5625

$ws.Range("H132").Value = 1787.5807
$ws.Range("I132").Value = 1473.5
$ws.Range("K132").Value = 4420.5
$ws.Range("M132").Value = -1890.5

$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -95100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 852.2727
$ws.Range("I22").Value = 504
$ws.Range("K22").Value = 504
$ws.Range("M22").Value = -209

$ws.Range("H27").Value = 852.2727
$ws.Range("I27").Value = 504
$ws.Range("K27").Value = 504
$ws.Range("M27").Value = -397

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 362.31818
$ws.Range("I107").Value = 281.82352
$ws.Range("J107").Value = 636
$ws.Range("K107").Value = 845.47056
$ws.Range("L107").Value = 1908
$ws.Range("M107").Value = 1074.52944
$ws.Range("N107").Value = -5748

$ws.Range("H132").Value = 1476
$ws.Range("I132").Value = 1194.4
$ws.Range("K132").Value = 3583.2
$ws.Range("M132").Value = -1053.2

$ws.Range("H136").Value = 436.4762
$ws.Range("I136").Value = 377.1579
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 1131.4737
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = 1418.5263
$ws.Range("N136").Value = -8100

$ws.Range("H137").Value = 31961.6
$ws.Range("J137").Value = 31961.6
$ws.Range("L137").Value = 31961.6
$ws.Range("N137").Value = -42161.6
